# Auto-generated edit script: updates crypto price/volume table per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.119.71'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.22%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.425.03'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.96%  '

$ws.Range('E4').Value = '  -0.22%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '561.17'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.55%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.63'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.93%  '

$ws.Range('E7').Value = '  +0.16%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.526'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.64%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.420.36'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.47%  '

$ws.Range('E10').Value = '  -6.10%  '

$ws.Range('E11').Value = '  +0.88%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.16'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.44%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.351'
$ws.Range('D13').Style = 'Normal'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.34'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.53%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000172'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.97%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.867.47'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.76%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.307.25'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.88%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.418.35'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.88%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.00'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.76%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.09'
$ws.Range('D20').Style = 'Normal'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '322.95'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.83%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.11'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.18%  '

$ws.Range('E23').Value = '  +2.89%  '

$ws.Range('E24').Value = '  +0.13%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '64.80'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.02%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '617.90'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.45%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.93'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.05%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.562.77'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.46%  '

$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.46%  '

$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0947'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -10.39%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.43'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.26%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.99'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.58%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.85'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.67%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.134'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -7.23%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.97'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.42%  '

$ws.Range('E36').Value = '  +0.21%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.44'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.62%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.373'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.40%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.59'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.03%  '

$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '146.58'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.47%  '

$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.18'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.41%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.72'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.31%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.50'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.60%  '

$ws.Range('E44').Value = '  +0.01%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.46'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -8.44%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '144.78'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.03%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.66'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.88%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0519'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.76%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.94'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.97%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.590'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.25%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0227'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.31%  '
